$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 76
$ws1.Range("F4").Value = 711
$ws1.Range("F5").Value = 132
$ws1.Range("F8").Value = 2695
$ws1.Range("F10").Value = 1702
$ws1.Range("F13").Value = 699
$ws1.Range("F14").Value = 860
$ws1.Range("F15").Value = 140
$ws1.Range("F17").Value = 1106
$ws1.Range("F21").Value = 6065
$ws1.Range("F22").Value = 243
$ws1.Range("F23").Value = 1233
$ws1.Range("F24").Value = 131
$ws1.Range("F25").Value = 170
$ws1.Range("F26").Value = 153
$ws1.Range("F27").Value = 282
$ws1.Range("F28").Value = 242
$ws1.Range("F30").Value = 1075
$ws1.Range("F31").Value = 869
$ws1.Range("F35").Value = 441
$ws1.Range("F36").Value = 1266
$ws1.Range("F38").Value = 132
$ws1.Range("F41").Value = 157

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 76
$ws4.Range("F4").Value = 711
$ws4.Range("F5").Value = 132
$ws4.Range("F11").Value = 2695
$ws4.Range("F13").Value = 1702
$ws4.Range("F16").Value = 699
$ws4.Range("F18").Value = 860
$ws4.Range("F19").Value = 140
$ws4.Range("F21").Value = 1106
$ws4.Range("F25").Value = 6065
$ws4.Range("F26").Value = 243
$ws4.Range("F27").Value = 1233
$ws4.Range("F28").Value = 131
$ws4.Range("F29").Value = 170
$ws4.Range("F30").Value = 153
$ws4.Range("F31").Value = 282
$ws4.Range("F32").Value = 242
$ws4.Range("F34").Value = 1075
$ws4.Range("F35").Value = 869
$ws4.Range("F39").Value = 441
$ws4.Range("F40").Value = 1266
$ws4.Range("F42").Value = 132
$ws4.Range("F45").Value = 157
